$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in header cell A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 09:52"

# Row 13: Rusia - updated stats
$ws.Range("B13").Value = 62773
$ws.Range("C13").Value = 4774
$ws.Range("D13").Value = 4891
$ws.Range("E13").Value = 57327
$ws.Range("G13").Value = 42
$ws.Range("H13").Value = 555

# Row 29: Ecuador -> Singapur
$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 11178
$ws.Range("C29").Value = 1037
$ws.Range("D29").Value = 896
$ws.Range("E29").Value = 10270
$ws.Range("F29").Value = 27
$ws.Range("H29").Value = 12

# Row 30: Corea del Sur -> Ecuador
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 10850
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1262
$ws.Range("E30").Value = 9051
$ws.Range("F30").Value = 141
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 537

# Row 31: Mexico -> Corea del Sur
$ws.Range("A31").Value = "Corea del Sur"
$ws.Range("B31").Value = 10702
$ws.Range("C31").Value = 8
$ws.Range("D31").Value = 8411
$ws.Range("E31").Value = 2051
$ws.Range("F31").Value = 55
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 240

# Row 32: Pakistan -> Mexico
$ws.Range("A32").Value = "Mexico"
$ws.Range("B32").Value = 10544
$ws.Range("C32").Value = 1043
$ws.Range("D32").Value = 2627
$ws.Range("E32").Value = 6947
$ws.Range("F32").Value = 378
$ws.Range("G32").Value = 113
$ws.Range("H32").Value = 970

# Row 33: Polonia -> Pakistan
$ws.Range("A33").Value = "Pakistan"
$ws.Range("B33").Value = 10513
$ws.Range("C33").Value = 437
$ws.Range("D33").Value = 2337
$ws.Range("E33").Value = 7952
$ws.Range("F33").Value = 60
$ws.Range("G33").Value = 12
$ws.Range("H33").Value = 224

# Row 34: Singapur -> Polonia
$ws.Range("A34").Value = "Polonia"
$ws.Range("B34").Value = 10169
$ws.Range("D34").Value = 1740
$ws.Range("E34").Value = 8003
$ws.Range("F34").Value = 160
$ws.Range("H34").Value = 426

# Row 69: Uzbekistan - updated stats
$ws.Range("D69").Value = 454
$ws.Range("E69").Value = 1255

# Row 72: Armenia - updated stats
$ws.Range("F72").Value = 10

# Row 78: Republica de Macedonia -> Eslovaquia
$ws.Range("A78").Value = "Eslovaquia"
$ws.Range("B78").Value = 1325
$ws.Range("C78").Value = 81
$ws.Range("D78").Value = 288
$ws.Range("E78").Value = 1022
$ws.Range("F78").Value = 9
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 15

# Row 79: Eslovaquia -> Republica de Macedonia
$ws.Range("A79").Value = "Republica de Macedonia"
$ws.Range("B79").Value = 1259
$ws.Range("D79").Value = 272
$ws.Range("E79").Value = 931
$ws.Range("H79").Value = 56

# Row 91: Letonia - updated stats
$ws.Range("B91").Value = 778
$ws.Range("C91").Value = 17
$ws.Range("E91").Value = 634
$ws.Range("F91").Value = 6

# Row 117: Montenegro - updated stats
$ws.Range("B117").Value = 316
$ws.Range("C117").Value = 1
$ws.Range("E117").Value = 195
